# issue #5: stock data from json to db
#
# The 股票 (stock) sheet gains three new columns:
#   - "category"    inserted right after "property_category" (col H), pushing
#                    the existing date / legislator_name / legislator_id
#                    columns one place to the right.
#   - "source_file" appended after the (now shifted) legislator_id column.
#   - "index"       appended after source_file, duplicating column A's id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 6

# --- Insert a new column at I (9), shifting date/legislator_name/
# legislator_id from I:K to J:L. Excel's native column Insert carries the
# data, shared-string typing and formatting along for the ride, so no
# explicit re-typing (and no risk of "2013-12-03" text being reinterpreted
# as a date serial) is needed.
$ws.Columns.Item(9).Insert()

# --- Header row ---
$ws.Cells.Item(1, 9).Value2 = "category"

# --- "category" data column: every row is "normal" ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value2 = "normal"
}

# --- Append "source_file" / "index" columns after legislator_id (L) ---
# Copy formatting from the existing legislator_id column (L) so the new
# M/N columns keep the workbook's existing header/data styles instead of
# picking up the blank default style.
$ws.Range("L1:L6").Copy() | Out-Null
$ws.Range("M1:N6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 13).Value2 = "source_file"
$ws.Cells.Item(1, 14).Value2 = "index"

for ($r = 2; $r -le $lastRow; $r++) {
    $idVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 13).Value2 = "tmp9b3d1"
    $ws.Cells.Item($r, 14).Value2 = $idVal
}
